$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11: Key Stage 4 (KS4) destinations - update source link and periods
$ws.Range("B11").Value = "<a href = 'https://explore-education-statistics.service.gov.uk/find-statistics/key-stage-4-destination-measures/'>Key stage 4 destination measures</a>"
$ws.Range("C11").Value = "Aug 2023 -  Jul 2024 (22/23 learners) (16/10/25)"
$ws.Range("D11").Value = "Aug 2023 -  Jul 2024 (22/23 learners) (Feb 26)"

# Row 12: Key Stage 5 (KS5) destinations - update source link and periods
$ws.Range("B12").Value = "<a href = 'https://explore-education-statistics.service.gov.uk/find-statistics/16-18-destination-measures'>16-18 destination measures</a>"
$ws.Range("C12").Value = "Aug 2023 -  Jul 2024 (22/23 learners) (16/10/25)"
$ws.Range("D12").Value = "Aug 2023 -  Jul 2024 (22/23 learners) (Feb 26)"

# Update selection to reflect the last-active cell recorded in the saved file
$ws.Range("D12").Select()
